$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 values are written as literal text (not auto-converted to numbers/dates)
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.657.35'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '2.125.39'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  +0.64%  '
$ws.Range("D5").Value = '352.41'
$ws.Range("E5").Value = '  +5.42%  '
$ws.Range("D6").Value = '1.010'
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("D8").Value = '0.4552'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '54.01'
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").Value = '0.09102'
$ws.Range("E10").Value = '  +1.80%  '
$ws.Range("D11").Value = '1.183'
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").Value = '24.63'
$ws.Range("D13").Value = '2.132.63'
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("D14").Value = '6.861'
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("D16").Value = '102.51'
$ws.Range("E16").Value = '  +6.07%  '
$ws.Range("D17").Value = '0.00001178'
$ws.Range("E17").Value = '  +2.85%  '
$ws.Range("D19").Value = '0.06711'
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("E20").Value = '  +1.20%  '
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").Value = '30.735.70'
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").Value = '12.90'
$ws.Range("E24").Value = '  +3.17%  '
$ws.Range("D25").Value = '2.396'
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("D26").Value = '2.386.72'
$ws.Range("E26").Value = '  +1.52%  '
$ws.Range("D27").Value = '22.47'
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("D29").Value = '164.74'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").Value = '136.42'
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("D31").Value = '1.201'
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("D33").Value = '1.671'
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").Value = '6.392'
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").Value = '4.019'
$ws.Range("E35").Value = '  +1.99%  '
$ws.Range("D36").Value = '6.192'
$ws.Range("E36").Value = '  +8.30%  '
$ws.Range("D37").Value = '10.33'
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("D38").Value = '0.02657'
$ws.Range("E38").Value = '  +2.86%  '
$ws.Range("D39").Value = '0.06886'
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").Value = '12.57'
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("D42").Value = '0.6931'
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("E43").Value = '  +2.17%  '
$ws.Range("D44").Value = '14.78'
$ws.Range("E44").Value = '  +4.89%  '
$ws.Range("D45").Value = '2.345'
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("D46").Value = '0.6466'
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("D47").Value = '3.761'
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("D48").Value = '0.00000000367'
$ws.Range("E48").Value = '  +6.25%  '
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").Value = '83.01'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").Value = '0.07300'
$ws.Range("E51").Value = '  +2.20%  '
# Restore original (default) cell style so no stray formatting is introduced
$dataRange.Style = "Normal"
